# Update res_bus vm_pu values for the 380 kV case (commit: "case with 380 kV done")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.0413424210947
$ws.Cells.Item(2, 4).Value = 1.040434543382687
$ws.Cells.Item(2, 5).Value = 1.045000411321775
$ws.Cells.Item(2, 6).Value = 1.047462803727005
$ws.Cells.Item(2, 9).Value = 1.035147007136494
$ws.Cells.Item(2, 10).Value = 1.046424145623197
$ws.Cells.Item(2, 11).Value = 1.043216902010085
$ws.Cells.Item(2, 12).Value = 1.047769880631462
$ws.Cells.Item(2, 13).Value = 1.050225372359725
$ws.Cells.Item(2, 14).Value = 1.047910187936509
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.043258014193539
$ws.Cells.Item(3, 4).Value = 1.041837158137755
$ws.Cells.Item(3, 5).Value = 1.046866709610554
$ws.Cells.Item(3, 6).Value = 1.049395990408469
$ws.Cells.Item(3, 9).Value = 1.035610101446189
$ws.Cells.Item(3, 10).Value = 1.047980919184381
$ws.Cells.Item(3, 11).Value = 1.044428382097392
$ws.Cells.Item(3, 12).Value = 1.049444783996867
$ws.Cells.Item(3, 13).Value = 1.051967503494144
$ws.Cells.Item(3, 14).Value = 1.049469172294714
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.044493181957053
$ws.Cells.Item(4, 4).Value = 1.04274084691629
$ws.Cells.Item(4, 5).Value = 1.048070388920144
$ws.Cells.Item(4, 6).Value = 1.050642860900118
$ws.Cells.Item(4, 9).Value = 1.035906563362
$ws.Cells.Item(4, 10).Value = 1.048983727451384
$ws.Cells.Item(4, 11).Value = 1.045207883381809
$ws.Cells.Item(4, 12).Value = 1.050524207251708
$ws.Cells.Item(4, 13).Value = 1.05309035070554
$ws.Cells.Item(4, 14).Value = 1.05047340466448
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.045011428042567
$ws.Cells.Item(5, 4).Value = 1.043119840344952
$ws.Cells.Item(5, 5).Value = 1.048575493512828
$ws.Cells.Item(5, 6).Value = 1.051166102560471
$ws.Cells.Item(5, 9).Value = 1.036030438434452
$ws.Cells.Item(5, 10).Value = 1.049404242614403
$ws.Cells.Item(5, 11).Value = 1.045534545269863
$ws.Cells.Item(5, 12).Value = 1.050976974766767
$ws.Cells.Item(5, 13).Value = 1.053561356211751
$ws.Cells.Item(5, 14).Value = 1.050894517007264
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.045098384774314
$ws.Cells.Item(6, 4).Value = 1.043183421664847
$ws.Cells.Item(6, 5).Value = 1.048660249329363
$ws.Cells.Item(6, 6).Value = 1.051253902498643
$ws.Cells.Item(6, 9).Value = 1.036051193368319
$ws.Cells.Item(6, 10).Value = 1.049474787021075
$ws.Cells.Item(6, 11).Value = 1.045589332648609
$ws.Cells.Item(6, 12).Value = 1.051052937060133
$ws.Cells.Item(6, 13).Value = 1.053640379749418
$ws.Cells.Item(6, 14).Value = 1.050965161595085
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.044500110757625
$ws.Cells.Item(7, 4).Value = 1.042745914628391
$ws.Cells.Item(7, 5).Value = 1.048077141747322
$ws.Cells.Item(7, 6).Value = 1.050649856154833
$ws.Cells.Item(7, 9).Value = 1.035908221554543
$ws.Cells.Item(7, 10).Value = 1.048989350556704
$ws.Cells.Item(7, 11).Value = 1.045212252322286
$ws.Cells.Item(7, 12).Value = 1.050530261142592
$ws.Cells.Item(7, 13).Value = 1.053096648360578
$ws.Cells.Item(7, 14).Value = 1.050479035755255
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041990719012413
$ws.Cells.Item(8, 4).Value = 1.040909379875173
$ws.Cells.Item(8, 5).Value = 1.045631964281979
$ws.Cells.Item(8, 6).Value = 1.048116981225743
$ws.Cells.Item(8, 9).Value = 1.035304176653915
$ws.Cells.Item(8, 10).Value = 1.046951213451327
$ws.Cells.Item(8, 11).Value = 1.04362724883188
$ws.Cells.Item(8, 12).Value = 1.048336834920993
$ws.Cells.Item(8, 13).Value = 1.050815063099652
$ws.Cells.Item(8, 14).Value = 1.04843800426141
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.037534513970916
$ws.Cells.Item(9, 4).Value = 1.037642622605152
$ws.Cells.Item(9, 5).Value = 1.041292112824263
$ws.Cells.Item(9, 6).Value = 1.043621857897874
$ws.Cells.Item(9, 9).Value = 1.034215041290589
$ws.Cells.Item(9, 10).Value = 1.043324230151333
$ws.Cells.Item(9, 11).Value = 1.04079985123122
$ws.Cells.Item(9, 12).Value = 1.044437523145791
$ws.Cells.Item(9, 13).Value = 1.046759771027764
$ws.Cells.Item(9, 14).Value = 1.04480587022911
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034539178568264
$ws.Cells.Item(10, 4).Value = 1.035443253282162
$ws.Cells.Item(10, 5).Value = 1.038376596872115
$ws.Cells.Item(10, 6).Value = 1.040602264265911
$ws.Cells.Item(10, 9).Value = 1.03347191202896
$ws.Cells.Item(10, 10).Value = 1.040881167648084
$ws.Cells.Item(10, 11).Value = 1.038890842227008
$ws.Cells.Item(10, 12).Value = 1.041813725580899
$ws.Cells.Item(10, 13).Value = 1.044031499680332
$ws.Cells.Item(10, 14).Value = 1.042359338296883
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033236022027254
$ws.Cells.Item(11, 4).Value = 1.034485572129359
$ws.Cells.Item(11, 5).Value = 1.037108558933066
$ws.Cells.Item(11, 6).Value = 1.039289010994653
$ws.Cells.Item(11, 9).Value = 1.03314599476213
$ws.Cells.Item(11, 10).Value = 1.039817080555004
$ws.Cells.Item(11, 11).Value = 1.038058299015407
$ws.Cells.Item(11, 12).Value = 1.040671557735082
$ws.Cells.Item(11, 13).Value = 1.042843964311055
$ws.Cells.Item(11, 14).Value = 1.041293740078079
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032751018568413
$ws.Cells.Item(12, 4).Value = 1.034129025086682
$ws.Cells.Item(12, 5).Value = 1.036636685569255
$ws.Cells.Item(12, 6).Value = 1.03880031896375
$ws.Cells.Item(12, 9).Value = 1.033024304891597
$ws.Cells.Item(12, 10).Value = 1.039420872831167
$ws.Cells.Item(12, 11).Value = 1.037748146225621
$ws.Cells.Item(12, 12).Value = 1.040246372960673
$ws.Cells.Item(12, 13).Value = 1.042401907350157
$ws.Cells.Item(12, 14).Value = 1.040896969693828
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032855097042575
$ws.Cells.Item(13, 4).Value = 1.034205543132152
$ws.Cells.Item(13, 5).Value = 1.036737943694212
$ws.Cells.Item(13, 6).Value = 1.038905185842379
$ws.Cells.Item(13, 9).Value = 1.033050436407738
$ws.Cells.Item(13, 10).Value = 1.039505904478202
$ws.Cells.Item(13, 11).Value = 1.037814716490523
$ws.Cells.Item(13, 12).Value = 1.040337619144154
$ws.Cells.Item(13, 13).Value = 1.042496773629387
$ws.Cells.Item(13, 14).Value = 1.040982122095554
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033195951097241
$ws.Cells.Item(14, 4).Value = 1.034456116713641
$ws.Cells.Item(14, 5).Value = 1.037069571585825
$ws.Cells.Item(14, 6).Value = 1.03924863390174
$ws.Cells.Item(14, 9).Value = 1.033135948732129
$ws.Cells.Item(14, 10).Value = 1.039784349587328
$ws.Cells.Item(14, 11).Value = 1.038032680342745
$ws.Cells.Item(14, 12).Value = 1.040636430999971
$ws.Cells.Item(14, 13).Value = 1.042807443344283
$ws.Cells.Item(14, 14).Value = 1.041260962628674
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033405835386244
$ws.Cells.Item(15, 4).Value = 1.034610393909215
$ws.Cells.Item(15, 5).Value = 1.037273782716218
$ws.Cells.Item(15, 6).Value = 1.03946012465086
$ws.Cells.Item(15, 9).Value = 1.033188552013739
$ws.Cells.Item(15, 10).Value = 1.039955781112411
$ws.Cells.Item(15, 11).Value = 1.038166854039472
$ws.Cells.Item(15, 12).Value = 1.040820414507262
$ws.Cells.Item(15, 13).Value = 1.04299873009545
$ws.Cells.Item(15, 14).Value = 1.041432637606186
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.034625532671801
$ws.Cells.Item(16, 4).Value = 1.035506697276132
$ws.Cells.Item(16, 5).Value = 1.038460632151875
$ws.Cells.Item(16, 6).Value = 1.04068929708898
$ws.Cells.Item(16, 9).Value = 1.033493454226828
$ws.Cells.Item(16, 10).Value = 1.040951654481889
$ws.Cells.Item(16, 11).Value = 1.038945968890111
$ws.Cells.Item(16, 12).Value = 1.041889398037681
$ws.Cells.Item(16, 13).Value = 1.044110180197168
$ws.Cells.Item(16, 14).Value = 1.042429925230077
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.035388948533446
$ws.Cells.Item(17, 4).Value = 1.036067481506731
$ws.Cells.Item(17, 5).Value = 1.039203593690895
$ws.Cells.Item(17, 6).Value = 1.04145876592832
$ws.Cells.Item(17, 9).Value = 1.033683597960508
$ws.Cells.Item(17, 10).Value = 1.041574656697491
$ws.Cells.Item(17, 11).Value = 1.039433086666867
$ws.Cells.Item(17, 12).Value = 1.042558307309125
$ws.Cells.Item(17, 13).Value = 1.044805692075131
$ws.Cells.Item(17, 14).Value = 1.043053812180285
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.035833643993724
$ws.Cells.Item(18, 4).Value = 1.036394063579374
$ws.Cells.Item(18, 5).Value = 1.039636412030817
$ws.Cells.Item(18, 6).Value = 1.041907031238863
$ws.Cells.Item(18, 9).Value = 1.033794106937875
$ws.Cells.Item(18, 10).Value = 1.041937444318903
$ws.Cells.Item(18, 11).Value = 1.039716643343127
$ws.Cells.Item(18, 12).Value = 1.042947889067682
$ws.Cells.Item(18, 13).Value = 1.045210778251904
$ws.Cells.Item(18, 14).Value = 1.043417115001731
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.035985174069553
$ws.Cells.Item(19, 4).Value = 1.036505333045592
$ws.Cells.Item(19, 5).Value = 1.039783901339365
$ws.Cells.Item(19, 6).Value = 1.04205978518038
$ws.Cells.Item(19, 9).Value = 1.03383172026901
$ws.Cells.Item(19, 10).Value = 1.042061044612739
$ws.Cells.Item(19, 11).Value = 1.039813232560115
$ws.Cells.Item(19, 12).Value = 1.043080628446101
$ws.Cells.Item(19, 13).Value = 1.045348802163186
$ws.Cells.Item(19, 14).Value = 1.043540890822162
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03530710268227
$ws.Cells.Item(20, 4).Value = 1.036007367958767
$ws.Cells.Item(20, 5).Value = 1.03912393684332
$ws.Cells.Item(20, 6).Value = 1.041376266585581
$ws.Cells.Item(20, 9).Value = 1.03366323862941
$ws.Cells.Item(20, 10).Value = 1.041507876583205
$ws.Cells.Item(20, 11).Value = 1.039380882720169
$ws.Cells.Item(20, 12).Value = 1.042486600018758
$ws.Cells.Item(20, 13).Value = 1.044731131871282
$ws.Cells.Item(20, 14).Value = 1.042986937230576
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033095604611918
$ws.Cells.Item(21, 4).Value = 1.034382351877218
$ws.Cells.Item(21, 5).Value = 1.036971939600337
$ws.Cells.Item(21, 6).Value = 1.039147521835599
$ws.Cells.Item(21, 9).Value = 1.033110784931945
$ws.Cells.Item(21, 10).Value = 1.039702381087098
$ws.Cells.Item(21, 11).Value = 1.0379685206668
$ws.Cells.Item(21, 12).Value = 1.040548464303288
$ws.Cells.Item(21, 13).Value = 1.04271598538379
$ws.Cells.Item(21, 14).Value = 1.041178877723772
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031699619434981
$ws.Cells.Item(22, 4).Value = 1.033355876302218
$ws.Cells.Item(22, 5).Value = 1.035613860498308
$ws.Cells.Item(22, 6).Value = 1.037741050874226
$ws.Cells.Item(22, 9).Value = 1.032759788096785
$ws.Cells.Item(22, 10).Value = 1.038561637675406
$ws.Cells.Item(22, 11).Value = 1.03707524313219
$ws.Cells.Item(22, 12).Value = 1.039324472313135
$ws.Cells.Item(22, 13).Value = 1.041443453485811
$ws.Cells.Item(22, 14).Value = 1.04003651432559
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032440191139284
$ws.Cells.Item(23, 4).Value = 1.033900488497388
$ws.Cells.Item(23, 5).Value = 1.036334289857767
$ws.Cells.Item(23, 6).Value = 1.038487147158313
$ws.Cells.Item(23, 9).Value = 1.032946206612704
$ws.Cells.Item(23, 10).Value = 1.039166901883751
$ws.Cells.Item(23, 11).Value = 1.037549292119813
$ws.Cells.Item(23, 12).Value = 1.039973854594223
$ws.Cells.Item(23, 13).Value = 1.042118579481723
$ws.Cells.Item(23, 14).Value = 1.040642638078536
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.035344087128091
$ws.Cells.Item(24, 4).Value = 1.036034532270896
$ws.Cells.Item(24, 5).Value = 1.039159932007302
$ws.Cells.Item(24, 6).Value = 1.04141354619487
$ws.Cells.Item(24, 9).Value = 1.033672439366356
$ws.Cells.Item(24, 10).Value = 1.041538053492631
$ws.Cells.Item(24, 11).Value = 1.039404473200074
$ws.Cells.Item(24, 12).Value = 1.04251900325339
$ws.Cells.Item(24, 13).Value = 1.04476482425308
$ws.Cells.Item(24, 14).Value = 1.043017156994675
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038690766469905
$ws.Cells.Item(25, 4).Value = 1.038490879064143
$ws.Cells.Item(25, 5).Value = 1.042417895824892
$ws.Cells.Item(25, 6).Value = 1.044787876679053
$ws.Cells.Item(25, 9).Value = 1.034499581602679
$ws.Cells.Item(25, 10).Value = 1.044266223519708
$ws.Cells.Item(25, 11).Value = 1.0415349756237
$ws.Cells.Item(25, 12).Value = 1.045449771162648
$ws.Cells.Item(25, 13).Value = 1.047812424960678
$ws.Cells.Item(25, 14).Value = 1.045749201336116
